# Apply cryptos list update (data refresh) to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.235.85"
$ws.Range("E2").Value = "  -0.37%  "

# Row 3
$ws.Range("D3").Value = "3.672.77"
$ws.Range("E3").Value = "  -0.44%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "673.85"
$ws.Range("E5").Value = "  -1.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.73"
$ws.Range("E6").Value = "  -3.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.493"
$ws.Range("E8").Value = "  -1.48%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.146"
$ws.Range("E9").Value = "  -1.79%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.97"
$ws.Range("E10").Value = "  -5.79%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.436"
$ws.Range("E11").Value = "  -2.44%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000231"
$ws.Range("E12").Value = "  -3.52%  "

# Row 13
$ws.Range("D13").Value = "4.292.42"
$ws.Range("E13").Value = "  -0.48%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.30"
$ws.Range("E14").Value = "  -4.00%  "

# Row 15
$ws.Range("D15").Value = "3.673.79"
$ws.Range("E15").Value = "  -0.34%  "

# Row 16
$ws.Range("D16").Value = "69.184.99"
$ws.Range("E16").Value = "  -0.50%  "

# Row 17
$ws.Range("E17").Value = "  +1.56%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.04"
$ws.Range("E18").Value = "  -1.47%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.42"
$ws.Range("E19").Value = "  -3.20%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.28"
$ws.Range("E20").Value = "  -3.45%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.97"
$ws.Range("E21").Value = "  +0.63%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.648"
$ws.Range("E22").Value = "  -2.91%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.66"
$ws.Range("E23").Value = "  -0.86%  "

# Row 24
$ws.Range("D24").Value = "3.817.18"

# Row 25
$ws.Range("E25").Value = "  -0.04%  "

# Row 26
$ws.Range("E26").Value = "  -7.38%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("E27").Value = "  -5.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.03"
$ws.Range("E28").Value = "  -5.64%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.67"
$ws.Range("E29").Value = "  -1.99%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.74"
$ws.Range("E30").Value = "  -6.04%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.63"
$ws.Range("E31").Value = "  -4.00%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.18%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.98"
$ws.Range("E33").Value = "  -5.22%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.83"
$ws.Range("E34").Value = "  -1.22%  "

# Row 35
$ws.Range("D35").Value = "3.665.88"
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
$ws.Range("E36").Value = "  -5.10%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.14"
$ws.Range("E37").Value = "  -4.59%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.23"
$ws.Range("E38").Value = "  -2.11%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.22"
$ws.Range("E41").Value = "  -1.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "175.14"
$ws.Range("E42").Value = "  +8.53%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0901"
$ws.Range("E43").Value = "  -4.23%  "

# Row 44
$ws.Range("E44").Value = "  -1.84%  "

# Row 45
$ws.Range("E45").Value = "  -1.68%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.70"
$ws.Range("E46").Value = "  -5.66%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.95"
$ws.Range("E47").Value = "  -8.62%  "

# Row 48
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000276"
$ws.Range("E48").Value = "  -4.72%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.28"
$ws.Range("E49").Value = "  -5.44%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.08"
$ws.Range("E50").Value = "  -3.90%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.78"
$ws.Range("E51").Value = "  -3.19%  "
